# revise: shorter presentation transitions
#
# The deck-wide "Glitter" slide transition (with a markup-compatibility
# fallback to a plain Fade) is removed from every slide. We also tighten
# up a split text run in the comparison table on slide 8 ("Emotive " +
# "verbal expression" -> "Emotive verbal expression").

$p = $ppt.ActivePresentation

# 1. Strip the slide transition (glitter/fade) from every slide by
#    clearing the entry effect and its duration - the COM equivalent of
#    picking "None" in the Transitions gallery for each slide.
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    $trans = $slide.SlideShowTransition
    $trans.EntryEffect = 0
    $trans.Duration = 0
}

# 2. Merge the two runs of the third bullet in the "Features" column of
#    the comparison table on slide 8 into a single run of text.
$slide8 = $p.Slides.Item(8)
$tableShape = $slide8.Shapes.Item(3)
$cell = $tableShape.Table.Cell(3, 7)
$paragraph = $cell.Shape.TextFrame.TextRange.Paragraphs(3)
$paragraph.Text = "Emotive verbal expression"
